$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36 - this shifts the existing rows 36:46 down to 37:47
# (keeping all of their data/styles intact) and makes room for a new weekly record.
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with this week's price record for
# Vega Monumental Concepción / Chirimoya / Cultivar IV Región / Primera.
$ws.Range("A36").Value2 = 11
$ws.Range("B36").Value2 = "Vega Monumental Concepción"
$ws.Range("C36").Value2 = "Bíobío"
$ws.Range("D36").Value2 = 44876
$ws.Range("E36").Value2 = 8
$ws.Range("F36").Value2 = "Fruta"
$ws.Range("G36").Value2 = 100107
$ws.Range("H36").Value2 = "Otros"
$ws.Range("I36").Value2 = 100107002
$ws.Range("J36").Value2 = "Chirimoya"
$ws.Range("K36").Value2 = "Cultivar IV Región"
$ws.Range("L36").Value2 = "Primera"
$ws.Range("M36").Value2 = 140
$ws.Range("N36").Value2 = 25000
$ws.Range("O36").Value2 = 26000
$ws.Range("P36").Value2 = 25429
$ws.Range("Q36").Value2 = "$/bandeja 10 kilos"
$ws.Range("R36").Value2 = "Provincia de Limarí"
$ws.Range("S36").Value2 = 2543
$ws.Range("T36").Value2 = 10
